$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'331.66"
$ws.Range("E2").Value = "'1.07%"
$ws.Range("E3").Value = "'1.60%"
$ws.Range("D4").Value = "'5.558"
$ws.Range("E4").Value = "'-0.37%"
$ws.Range("D5").Value = "'0.08262"
$ws.Range("E5").Value = "'2.80%"
$ws.Range("D6").Value = "'2.042"
$ws.Range("E6").Value = "'3.50%"
$ws.Range("D7").Value = "'0.9771"
$ws.Range("E7").Value = "'3.12%"
$ws.Range("D8").Value = "'0.1122"
$ws.Range("E8").Value = "'-4.18%"
$ws.Range("E9").Value = "'2.96%"
$ws.Range("D10").Value = "'10.29"
$ws.Range("E10").Value = "'-12.97%"
$ws.Range("D11").Value = "'0.1002"
$ws.Range("E11").Value = "'2.14%"
$ws.Range("D12").Value = "'0.04669"
$ws.Range("E12").Value = "'-1.13%"
$ws.Range("D13").Value = "'0.1059"
$ws.Range("E13").Value = "'-0.63%"
$ws.Range("D14").Value = "'0.001284"
$ws.Range("E14").Value = "'-0.23%"
$ws.Range("D15").Value = "'0.04106"
$ws.Range("E15").Value = "'-2.53%"
$ws.Range("D16").Value = "'0.005946"
$ws.Range("E16").Value = "'-0.51%"
$ws.Range("D17").Value = "'3.363"
$ws.Range("E17").Value = "'-0.20%"
$ws.Range("D18").Value = "'4.442"
$ws.Range("E18").Value = "'2.67%"
$ws.Range("D19").Value = "'2.644"
$ws.Range("E19").Value = "'3.62%"
$ws.Range("D20").Value = "'0.3351"
$ws.Range("E20").Value = "'-3.56%"
$ws.Range("D21").Value = "'0.1384"
$ws.Range("E21").Value = "'-1.74%"
$ws.Range("E22").Value = "'-0.75%"
$ws.Range("E23").Value = "'3.85%"
$ws.Range("D24").Value = "'0.004397"
$ws.Range("E24").Value = "'2.01%"
$ws.Range("D25").Value = "'0.0001282"
$ws.Range("E25").Value = "'7.48%"
$ws.Range("D26").Value = "'0.0003745"
$ws.Range("E26").Value = "'-0.19%"
$ws.Range("D38").Value = "'0.02784"
$ws.Range("E38").Value = "'7.37%"
$ws.Range("D39").Value = "'0.05746"
$ws.Range("E39").Value = "'4.36%"
$ws.Range("D40").Value = "'0.007650"
$ws.Range("E40").Value = "'1.34%"
$ws.Range("D41").Value = "'0.1424"
$ws.Range("E41").Value = "'1.57%"
$ws.Range("D42").Value = "'0.007562"
$ws.Range("E42").Value = "'-1.18%"
$ws.Range("D43").Value = "'0.001976"
$ws.Range("E43").Value = "'-2.17%"
$ws.Range("D44").Value = "'0.008314"
$ws.Range("E44").Value = "'-0.64%"
$ws.Range("E45").Value = "'-0.94%"
$ws.Range("D46").Value = "'0.00000000751"
$ws.Range("E46").Value = "'-0.08%"
$ws.Range("D47").Value = "'0.0005801"
$ws.Range("E47").Value = "'-0.18%"
$ws.Range("D48").Value = "'0.002524"
$ws.Range("E48").Value = "'9.64%"
$ws.Range("D49").Value = "'0.003788"
$ws.Range("E49").Value = "'-21.76%"
$ws.Range("D50").Value = "'0.00002103"
$ws.Range("E50").Value = "'-0.08%"
$ws.Range("D51").Value = "'0.0002003"
$ws.Range("E51").Value = "'-0.08%"
